# Update New Orleans xlsx: add a "State" column to hotel_info and
# reorder the sheet tabs so review_info comes before hotel_info.

$wb = $excel.ActiveWorkbook

$hotelSheet = $wb.Worksheets.Item("hotel_info")
$reviewSheet = $wb.Worksheets.Item("review_info")

# Insert a new "State" column between "Hotel_Name" (B) and "City" (C).
$hotelSheet.Range("C1:C2").EntireColumn.Insert()
$hotelSheet.Range("C1").Value = "State"
$hotelSheet.Range("C2").Value = "Louisiana"

# Reorder the sheet tabs: review_info first, hotel_info second.
$reviewSheet.Move($hotelSheet)
